$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("E2").Value = 0.000246323
$ws.Range("F2").Value = 0.016120654
$ws.Range("G2").Value = 0.0004682440455

$ws.Range("E3").Value = 0.009391366
$ws.Range("F3").Value = 0.015660926
$ws.Range("G3").Value = 0.01053818123628692
